# Workbook / sheet handles
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1. Add the new "metadata" sheet right after "data" -------------------
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# --- 2. Populate the metadata header row (row 1, columns B:G) -------------
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# --- 3. Populate the metadata data row (row 2) -----------------------------
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Choanal atresia"
$meta.Range("C2").Value = 3498
# Force "1.2" to stay text (not be auto-parsed as the number 1.2), the same
# way Excel keeps a leading-apostrophe entry as text; ClearFormats drops the
# visible "text-quote" cell style the apostrophe trick adds, leaving the cell
# plain (like its siblings B2/C2/E2/F2/G2) but still storing a text value.
$meta.Range("D2").Value = "'1.2"
$meta.Range("D2").ClearFormats()
$meta.Range("E2").Value = "2021-03-27T00:30:58.774320Z"
$meta.Range("F2").Value = "2021-10-05 14:33:26.241905"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3498/?format=json"

# --- 4. Match the formatting used by the "data" sheet's header/index cells -
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 5. Refresh the "time_taken" column on the "data" sheet ---------------
$data.Range("F2").Value = "2021-10-05 14:33:26.244474"
$data.Range("F3").Value = "2021-10-05 14:33:26.244482"
$data.Range("F4").Value = "2021-10-05 14:33:26.244484"
$data.Range("F5").Value = "2021-10-05 14:33:26.244486"
$data.Range("F6").Value = "2021-10-05 14:33:26.244489"
$data.Range("F7").Value = "2021-10-05 14:33:26.244491"
$data.Range("F8").Value = "2021-10-05 14:33:26.244494"
$data.Range("F9").Value = "2021-10-05 14:33:26.244496"
$data.Range("F10").Value = "2021-10-05 14:33:26.244498"
$data.Range("F11").Value = "2021-10-05 14:33:26.244500"
$data.Range("F12").Value = "2021-10-05 14:33:26.244503"
$data.Range("F13").Value = "2021-10-05 14:33:26.244505"
$data.Range("F14").Value = "2021-10-05 14:33:26.244507"
$data.Range("F15").Value = "2021-10-05 14:33:26.244509"
$data.Range("F16").Value = "2021-10-05 14:33:26.244511"
$data.Range("F17").Value = "2021-10-05 14:33:26.244513"

# --- 6. Re-select the "data" sheet so it stays the active tab -------------
$data.Activate()
